$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 93
$ws.Range("A93").Value = 91
$ws.Range("B93").Value = 6236251
$ws.Range("C93").Value = "Venezuela Primera Division"
$ws.Range("D93").Value = 45199.6875
$ws.Range("E93").Value = "Angostura FC"
$ws.Range("F93").Value = "Portuguesa"
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 2
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 2
$ws.Range("K93").Value = "A"
$ws.Range("L93").Value = 3.1
$ws.Range("M93").Value = 3.2
$ws.Range("N93").Value = 2.15
$ws.Range("O93").Value = 4
$ws.Range("P93").Value = 3.6
$ws.Range("Q93").Value = 1.75
$ws.Range("R93").Value = 0.75
$ws.Range("S93").Value = 1.8
$ws.Range("T93").Value = 2
$ws.Range("U93").Value = 2.5
$ws.Range("V93").Value = 1.95
$ws.Range("W93").Value = 1.85
$ws.Range("X93").Value = -1
$ws.Range("Y93").Value = -1
$ws.Range("Z93").Value = 0.75
$ws.Range("AA93").Value = -0.5
$ws.Range("AB93").Value = 0.5
$ws.Range("AC93").Value = 0.95
$ws.Range("AD93").Value = -1

# Row 94
$ws.Range("A94").Value = 92
$ws.Range("B94").Value = 6236252
$ws.Range("C94").Value = "Venezuela Primera Division"
$ws.Range("D94").Value = 45199.6875
$ws.Range("E94").Value = "Deportivo Tachira"
$ws.Range("F94").Value = "CD Hermanos Colmenares"
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = "H"
$ws.Range("L94").Value = 1.363
$ws.Range("M94").Value = 4.2
$ws.Range("N94").Value = 7.5
$ws.Range("O94").Value = 1.333
$ws.Range("P94").Value = 4.5
$ws.Range("Q94").Value = 8
$ws.Range("R94").Value = -1.5
$ws.Range("S94").Value = 2
$ws.Range("T94").Value = 1.8
$ws.Range("U94").Value = 2.5
$ws.Range("V94").Value = 1.925
$ws.Range("W94").Value = 1.875
$ws.Range("X94").Value = 0.333
$ws.Range("Y94").Value = -1
$ws.Range("Z94").Value = -1
$ws.Range("AA94").Value = -1
$ws.Range("AB94").Value = 0.8
$ws.Range("AC94").Value = -1
$ws.Range("AD94").Value = 0.875

# Row 95
$ws.Range("A95").Value = 93
$ws.Range("B95").Value = 6236254
$ws.Range("C95").Value = "Venezuela Primera Division"
$ws.Range("D95").Value = 45199.6875
$ws.Range("E95").Value = "Academia Puerto Cabello"
$ws.Range("F95").Value = "Estudiantes Merida"
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 1
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = "H"
$ws.Range("L95").Value = 1.727
$ws.Range("M95").Value = 3.4
$ws.Range("N95").Value = 4.333
$ws.Range("O95").Value = 1.666
$ws.Range("P95").Value = 3.4
$ws.Range("Q95").Value = 4.75
$ws.Range("R95").Value = -0.75
$ws.Range("S95").Value = 1.875
$ws.Range("T95").Value = 1.925
$ws.Range("U95").Value = 2.5
$ws.Range("V95").Value = 1.9
$ws.Range("W95").Value = 1.9
$ws.Range("X95").Value = 0.6659999999999999
$ws.Range("Y95").Value = -1
$ws.Range("Z95").Value = -1
$ws.Range("AA95").Value = 0.4375
$ws.Range("AB95").Value = -0.5
$ws.Range("AC95").Value = -1
$ws.Range("AD95").Value = 0.8999999999999999

# Row 102
$ws.Range("A102").Value = 100
$ws.Range("B102").Value = 6236616
$ws.Range("C102").Value = "Venezuela Primera Division"
$ws.Range("D102").Value = 45206.6875
$ws.Range("E102").Value = "UCV"
$ws.Range("F102").Value = "Metropolitanos FC"
$ws.Range("G102").Value = 3
$ws.Range("H102").Value = 2
$ws.Range("I102").Value = 3
$ws.Range("J102").Value = 1
$ws.Range("K102").Value = "H"
$ws.Range("L102").Value = 3.3
$ws.Range("M102").Value = 3.2
$ws.Range("N102").Value = 2.05
$ws.Range("O102").Value = 2.75
$ws.Range("P102").Value = 3.2
$ws.Range("Q102").Value = 2.3
$ws.Range("R102").Value = 0.25
$ws.Range("S102").Value = 1.75
$ws.Range("T102").Value = 2.05
$ws.Range("U102").Value = 2.5
$ws.Range("V102").Value = 1.975
$ws.Range("W102").Value = 1.825
$ws.Range("X102").Value = 1.75
$ws.Range("Y102").Value = -1
$ws.Range("Z102").Value = -1
$ws.Range("AA102").Value = 0.75
$ws.Range("AB102").Value = -1
$ws.Range("AC102").Value = 0.9750000000000001
$ws.Range("AD102").Value = -1

# Row 103
$ws.Range("A103").Value = 101
$ws.Range("B103").Value = 6236615
$ws.Range("C103").Value = "Venezuela Primera Division"
$ws.Range("D103").Value = 45206.6875
$ws.Range("E103").Value = "Deportivo Rayo Zuliano"
$ws.Range("F103").Value = "Academia Puerto Cabello"
$ws.Range("G103").Value = 1
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 1
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = "H"
$ws.Range("L103").Value = 2.375
$ws.Range("M103").Value = 3.3
$ws.Range("N103").Value = 2.625
$ws.Range("O103").Value = 2.45
$ws.Range("P103").Value = 3.2
$ws.Range("Q103").Value = 2.55
$ws.Range("R103").Value = 0
$ws.Range("S103").Value = 1.875
$ws.Range("T103").Value = 1.925
$ws.Range("U103").Value = 2.5
$ws.Range("V103").Value = 2
$ws.Range("W103").Value = 1.8
$ws.Range("X103").Value = 1.45
$ws.Range("Y103").Value = -1
$ws.Range("Z103").Value = -1
$ws.Range("AA103").Value = 0.875
$ws.Range("AB103").Value = -1
$ws.Range("AC103").Value = -1
$ws.Range("AD103").Value = 0.8

# Row 157
$ws.Range("A157").Value = 155
$ws.Range("B157").Value = 7920997
$ws.Range("C157").Value = "Venezuela Primera Division"
$ws.Range("D157").Value = 45360.79166666666
$ws.Range("E157").Value = "Carabobo"
$ws.Range("F157").Value = "UCV"
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 1
$ws.Range("I157").Value = 0
$ws.Range("J157").Value = 0
$ws.Range("K157").Value = "A"
$ws.Range("L157").Value = 1.833
$ws.Range("M157").Value = 3.1
$ws.Range("N157").Value = 4.2
$ws.Range("O157").Value = 1.833
$ws.Range("P157").Value = 3.1
$ws.Range("Q157").Value = 4.2
$ws.Range("R157").Value = -0.5
$ws.Range("S157").Value = 1.9
$ws.Range("T157").Value = 1.9
$ws.Range("U157").Value = 2
$ws.Range("V157").Value = 1.85
$ws.Range("W157").Value = 1.95
$ws.Range("X157").Value = -1
$ws.Range("Y157").Value = -1
$ws.Range("Z157").Value = 3.2
$ws.Range("AA157").Value = -1
$ws.Range("AB157").Value = 0.8999999999999999
$ws.Range("AC157").Value = -1
$ws.Range("AD157").Value = 0.95

# Row 158
$ws.Range("A158").Value = 156
$ws.Range("B158").Value = 7920998
$ws.Range("C158").Value = "Venezuela Primera Division"
$ws.Range("D158").Value = 45360.79166666666
$ws.Range("E158").Value = "Zamora"
$ws.Range("F158").Value = "Caracas"
$ws.Range("G158").Value = 2
$ws.Range("H158").Value = 2
$ws.Range("I158").Value = 1
$ws.Range("J158").Value = 0
$ws.Range("K158").Value = "D"
$ws.Range("L158").Value = 3.75
$ws.Range("M158").Value = 3.2
$ws.Range("N158").Value = 1.909
$ws.Range("O158").Value = 3
$ws.Range("P158").Value = 2.9
$ws.Range("Q158").Value = 2.375
$ws.Range("R158").Value = 0.25
$ws.Range("S158").Value = 1.8
$ws.Range("T158").Value = 2
$ws.Range("U158").Value = 2
$ws.Range("V158").Value = 1.825
$ws.Range("W158").Value = 1.975
$ws.Range("X158").Value = -1
$ws.Range("Y158").Value = 1.9
$ws.Range("Z158").Value = -1
$ws.Range("AA158").Value = 0.4
$ws.Range("AB158").Value = -0.5
$ws.Range("AC158").Value = 0.825
$ws.Range("AD158").Value = -1

# Row 162
$ws.Range("A162").Value = 160
$ws.Range("B162").Value = 7952905
$ws.Range("C162").Value = "Venezuela Primera Division"
$ws.Range("D162").Value = 45366.83333333334
$ws.Range("E162").Value = "Angostura FC"
$ws.Range("F162").Value = "Deportivo Tachira"
$ws.Range("G162").Value = 2
$ws.Range("H162").Value = 0
$ws.Range("I162").Value = 0
$ws.Range("J162").Value = 0
$ws.Range("K162").Value = "H"
$ws.Range("L162").Value = 3.6
$ws.Range("M162").Value = 3.6
$ws.Range("N162").Value = 1.8
$ws.Range("O162").Value = 3.75
$ws.Range("P162").Value = 2.875
$ws.Range("Q162").Value = 2.1
$ws.Range("R162").Value = 0.25
$ws.Range("S162").Value = 1.95
$ws.Range("T162").Value = 1.85
$ws.Range("U162").Value = 2
$ws.Range("V162").Value = 2.025
$ws.Range("W162").Value = 1.775
$ws.Range("X162").Value = 2.75
$ws.Range("Y162").Value = -1
$ws.Range("Z162").Value = -1
$ws.Range("AA162").Value = 0.95
$ws.Range("AB162").Value = -1
$ws.Range("AC162").Value = 0
$ws.Range("AD162").Value = 0

# Row 163
$ws.Range("A163").Value = 161
$ws.Range("B163").Value = 7952893
$ws.Range("C163").Value = "Venezuela Primera Division"
$ws.Range("D163").Value = 45366.83333333334
$ws.Range("E163").Value = "UCV"
$ws.Range("F163").Value = "Deportivo La Guaira"
$ws.Range("G163").Value = 1
$ws.Range("H163").Value = 1
$ws.Range("I163").Value = 1
$ws.Range("J163").Value = 1
$ws.Range("K163").Value = "D"
$ws.Range("L163").Value = 2.1
$ws.Range("M163").Value = 3
$ws.Range("N163").Value = 3.25
$ws.Range("O163").Value = 2.25
$ws.Range("P163").Value = 3.1
$ws.Range("Q163").Value = 2.9
$ws.Range("R163").Value = -0.25
$ws.Range("S163").Value = 2.025
$ws.Range("T163").Value = 1.775
$ws.Range("U163").Value = 2
$ws.Range("V163").Value = 1.8
$ws.Range("W163").Value = 2
$ws.Range("X163").Value = -1
$ws.Range("Y163").Value = 2.1
$ws.Range("Z163").Value = -1
$ws.Range("AA163").Value = -0.5
$ws.Range("AB163").Value = 0.3875
$ws.Range("AC163").Value = 0
$ws.Range("AD163").Value = 0

# Row 173
$ws.Range("A173").Value = 171
$ws.Range("B173").Value = 7958193
$ws.Range("C173").Value = "Venezuela Primera Division"
$ws.Range("D173").Value = 45371.89583333334
$ws.Range("E173").Value = "Zamora"
$ws.Range("F173").Value = "Academia Puerto Cabello"
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 0
$ws.Range("I173").Value = 0
$ws.Range("J173").Value = 0
$ws.Range("K173").Value = "D"
$ws.Range("L173").Value = 3.75
$ws.Range("M173").Value = 3.3
$ws.Range("N173").Value = 1.85
$ws.Range("O173").Value = 3.1
$ws.Range("P173").Value = 3.2
$ws.Range("Q173").Value = 2.1
$ws.Range("R173").Value = 0.25
$ws.Range("S173").Value = 1.875
$ws.Range("T173").Value = 1.925
$ws.Range("U173").Value = 2.25
$ws.Range("V173").Value = 2.025
$ws.Range("W173").Value = 1.775
$ws.Range("X173").Value = -1
$ws.Range("Y173").Value = 2.2
$ws.Range("Z173").Value = -1
$ws.Range("AA173").Value = 0.4375
$ws.Range("AB173").Value = -0.5
$ws.Range("AC173").Value = -1
$ws.Range("AD173").Value = 0.7749999999999999

# Row 174
$ws.Range("A174").Value = 172
$ws.Range("B174").Value = 7958192
$ws.Range("C174").Value = "Venezuela Primera Division"
$ws.Range("D174").Value = 45371.89583333334
$ws.Range("E174").Value = "Deportivo Tachira"
$ws.Range("F174").Value = "Monagas"
$ws.Range("G174").Value = 1
$ws.Range("H174").Value = 0
$ws.Range("I174").Value = 0
$ws.Range("J174").Value = 0
$ws.Range("K174").Value = "H"
$ws.Range("L174").Value = 1.666
$ws.Range("M174").Value = 3.4
$ws.Range("N174").Value = 4.5
$ws.Range("O174").Value = 1.95
$ws.Range("P174").Value = 3.25
$ws.Range("Q174").Value = 3.5
$ws.Range("R174").Value = -0.5
$ws.Range("S174").Value = 1.975
$ws.Range("T174").Value = 1.825
$ws.Range("U174").Value = 2.25
$ws.Range("V174").Value = 2.025
$ws.Range("W174").Value = 1.775
$ws.Range("X174").Value = 0.95
$ws.Range("Y174").Value = -1
$ws.Range("Z174").Value = -1
$ws.Range("AA174").Value = 0.9750000000000001
$ws.Range("AB174").Value = -1
$ws.Range("AC174").Value = -1
$ws.Range("AD174").Value = 0.7749999999999999

# New row 233 (append)
$ws.Range("A232").Copy() | Out-Null
$ws.Range("A233").PasteSpecial(-4122) | Out-Null
$ws.Range("D232").Copy() | Out-Null
$ws.Range("D233").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A233").Value = 231
$ws.Range("B233").Value = 8145483
$ws.Range("C233").Value = "Venezuela Primera Division"
$ws.Range("D233").Value = 45438.875
$ws.Range("E233").Value = "Portuguesa"
$ws.Range("F233").Value = "Carabobo"
$ws.Range("G233").Value = 2
$ws.Range("H233").Value = 1
$ws.Range("K233").Value = "H"
$ws.Range("L233").Value = 2.3
$ws.Range("M233").Value = 2.9
$ws.Range("N233").Value = 3.1
$ws.Range("O233").Value = 1.8
$ws.Range("P233").Value = 3.3
$ws.Range("Q233").Value = 4.2
$ws.Range("R233").Value = -0.5
$ws.Range("S233").Value = 1.8
$ws.Range("T233").Value = 2
$ws.Range("U233").Value = 2
$ws.Range("V233").Value = 1.75
$ws.Range("W233").Value = 2.05
$ws.Range("X233").Value = 0.8
$ws.Range("Y233").Value = -1
$ws.Range("Z233").Value = -1
$ws.Range("AA233").Value = 0.8
$ws.Range("AB233").Value = -1
$ws.Range("AC233").Value = 0.75
$ws.Range("AD233").Value = -1
